$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume number and report date range) ---
$ws.Range("A8").Value = "Volume 30   Number  17"
$ws.Range("C9").Value = "Report Covering the Week  4/24/2023  Through  4/30/2023"

# --- Number formats used by each style bucket, keyed by the cellXfs index they map to ---
$numFmts = @{
    15 = '#,##0.0;"-"#,##0.0'
    16 = "#,##0"
    18 = "#,##0"
    19 = '#,##0.00;"-"#,##0.00'
}

# --- Data-table driven cell updates for the CompStat grid (rows 15-29) ---
$cellUpdates = @(
    @{addr="C15"; val=2; style=16},
    @{addr="F15"; val=3; style=16},
    @{addr="I15"; val=8; style=16},
    @{addr="K15"; val=33.333333333333; style=15},
    @{addr="L15"; val=60; style=15},
    @{addr="M15"; val=-20; style=15},
    @{addr="N15"; val=-20; style=15},
    @{addr="C16"; val=5; style=16},
    @{addr="D16"; val=4; style=16},
    @{addr="E16"; val=25; style=15},
    @{addr="F16"; val=16; style=16},
    @{addr="G16"; val=14; style=16},
    @{addr="H16"; val=14.285714285714; style=15},
    @{addr="I16"; val=70; style=16},
    @{addr="J16"; val=68; style=16},
    @{addr="K16"; val=2.941176470588; style=15},
    @{addr="L16"; val=34.615384615384; style=15},
    @{addr="M16"; val=-14.634146341463; style=15},
    @{addr="N16"; val=-77.198697068403; style=15},
    @{addr="C17"; val=6; style=16},
    @{addr="D17"; val=6; style=16},
    @{addr="E17"; val=0; style=15},
    @{addr="F17"; val=28; style=16},
    @{addr="G17"; val=24; style=16},
    @{addr="H17"; val=16.666666666666; style=15},
    @{addr="I17"; val=103; style=16},
    @{addr="J17"; val=92; style=16},
    @{addr="K17"; val=11.95652173913; style=15},
    @{addr="L17"; val=24.096385542168; style=15},
    @{addr="M17"; val=157.5; style=15},
    @{addr="N17"; val=-15.573770491803; style=15},
    @{addr="C18"; val=3; style=16},
    @{addr="D18"; val=1; style=16},
    @{addr="E18"; val=200; style=15},
    @{addr="F18"; val=7; style=16},
    @{addr="G18"; val=7; style=16},
    @{addr="H18"; val=0; style=15},
    @{addr="I18"; val=33; style=16},
    @{addr="J18"; val=43; style=16},
    @{addr="K18"; val=-23.255813953488; style=15},
    @{addr="L18"; val=10; style=15},
    @{addr="M18"; val=-60.714285714285; style=15},
    @{addr="N18"; val=-91.624365482233; style=15},
    @{addr="C19"; val=15; style=16},
    @{addr="D19"; val=15; style=16},
    @{addr="E19"; val=0; style=15},
    @{addr="G19"; val=52; style=16},
    @{addr="H19"; val=0; style=15},
    @{addr="I19"; val=192; style=16},
    @{addr="J19"; val=244; style=16},
    @{addr="K19"; val=-21.311475409836; style=15},
    @{addr="L19"; val=34.265734265734; style=15},
    @{addr="M19"; val=84.615384615384; style=15},
    @{addr="N19"; val=6.077348066298; style=15},
    @{addr="D20"; val=2; style=16},
    @{addr="E20"; val=0; style=15},
    @{addr="F20"; val=7; style=16},
    @{addr="G20"; val=12; style=16},
    @{addr="H20"; val=-41.666666666666; style=15},
    @{addr="I20"; val=75; style=16},
    @{addr="J20"; val=95; style=16},
    @{addr="K20"; val=-21.052631578947; style=15},
    @{addr="L20"; val=70.454545454545; style=15},
    @{addr="M20"; val=-20.212765957446; style=15},
    @{addr="N20"; val=-93.438320209973; style=15},
    @{addr="C21"; val=33; style=18},
    @{addr="D21"; val=28; style=18},
    @{addr="E21"; val=17.857142857142; style=19},
    @{addr="F21"; val=113; style=18},
    @{addr="G21"; val=109; style=18},
    @{addr="H21"; val=3.669724770642; style=19},
    @{addr="I21"; val=481; style=18},
    @{addr="J21"; val=548; style=18},
    @{addr="K21"; val=-12.226277372262; style=19},
    @{addr="L21"; val=33.983286908078; style=19},
    @{addr="M21"; val=16.183574879227; style=19},
    @{addr="N21"; val=-77.823881973259; style=19},
    @{addr="C22"; val=2; style=16},
    @{addr="F22"; val=2; style=16},
    @{addr="I22"; val=7; style=16},
    @{addr="K22"; val=16.666666666666; style=15},
    @{addr="L22"; val=75; style=15},
    @{addr="M22"; val=-12.5; style=15},
    @{addr="C24"; val=23; style=16},
    @{addr="D24"; val=27; style=16},
    @{addr="E24"; val=-14.814814814814; style=15},
    @{addr="F24"; val=74; style=16},
    @{addr="G24"; val=120; style=16},
    @{addr="H24"; val=-38.333333333333; style=15},
    @{addr="I24"; val=418; style=16},
    @{addr="J24"; val=454; style=16},
    @{addr="K24"; val=-7.929515418502; style=15},
    @{addr="L24"; val=65.217391304347; style=15},
    @{addr="M24"; val=102.912621359223; style=15},
    @{addr="D25"; val=9; style=16},
    @{addr="E25"; val=-11.111111111111; style=15},
    @{addr="F25"; val=41; style=16},
    @{addr="G25"; val=37; style=16},
    @{addr="H25"; val=10.81081081081; style=15},
    @{addr="I25"; val=166; style=16},
    @{addr="J25"; val=132; style=16},
    @{addr="K25"; val=25.757575757575; style=15},
    @{addr="L25"; val=34.959349593495; style=15},
    @{addr="M25"; val=3.75; style=15},
    @{addr="C26"; val=2; style=16},
    @{addr="F26"; val=4; style=16},
    @{addr="H26"; val=300; style=15},
    @{addr="I26"; val=12; style=16},
    @{addr="K26"; val=0; style=15},
    @{addr="L26"; val=50; style=15},
    @{addr="D27"; val="STR0"; style=14},
    @{addr="E27"; val="STR***"; style=14},
    @{addr="F27"; val=3; style=16},
    @{addr="G27"; val=2; style=16},
    @{addr="H27"; val=50; style=15},
    @{addr="I27"; val=21; style=16},
    @{addr="K27"; val=10.526315789473; style=15},
    @{addr="L27"; val=50; style=15},
    @{addr="D28"; val=1; style=16},
    @{addr="E28"; val=-100; style=15},
    @{addr="G28"; val=1; style=16},
    @{addr="H28"; val=-100; style=15},
    @{addr="J28"; val=3; style=16},
    @{addr="K28"; val=33.333333333333; style=15},
    @{addr="L28"; val=33.333333333333; style=15},
    @{addr="D29"; val=1; style=16},
    @{addr="E29"; val=-100; style=15},
    @{addr="G29"; val=1; style=16},
    @{addr="H29"; val=-100; style=15},
    @{addr="J29"; val=3; style=16},
    @{addr="K29"; val=-33.333333333333; style=15},
    @{addr="L29"; val=0; style=15}
)

foreach ($u in $cellUpdates) {
    $addr = $u.addr
    $style = $u.style
    $val = $u.val

    if ($val -eq "STR0" -or $val -eq "STR***") {
        # Numeric/placeholder cell reverting to the sheet's standard text placeholder
        # ("0" or "***.*") — match it via a known reference cell that already carries
        # that placeholder with the exact target style, then overwrite the value.
        if ($val -eq "STR0") {
            $refAddr = "C23"
            $text = "0"
        } else {
            $refAddr = "E23"
            $text = "***.*"
        }
        $ws.Range($addr).NumberFormat = "@"
        $ws.Range($addr).Value = $text
        $ws.Range($refAddr).Copy()
        $ws.Range($addr).PasteSpecial(-4122)
    } else {
        $ws.Range($addr).Value = $val
        if ($numFmts.ContainsKey($style)) {
            $ws.Range($addr).NumberFormat = $numFmts[$style]
        }
    }
}

Write-Output "Applied $($cellUpdates.Count) cell updates plus header text changes."
